$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.397.16"
$ws.Range("E2").Value = "  -3.27%  "
$ws.Range("D3").Value = "2.942.74"
$ws.Range("E3").Value = "  -4.23%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'494.44"
$ws.Range("E5").Value = "  -6.65%  "
$ws.Range("D6").Value = "'133.36"
$ws.Range("E6").Value = "  -7.31%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("D8").Value = "'0.423"
$ws.Range("E8").Value = "  -5.91%  "
$ws.Range("E9").Value = "  -6.85%  "
$ws.Range("E10").Value = "  -7.61%  "
$ws.Range("D11").Value = "'0.351"
$ws.Range("E11").Value = "  -5.81%  "
$ws.Range("D12").Value = "3.446.35"
$ws.Range("E12").Value = "  -4.13%  "
$ws.Range("E13").Value = "  -3.73%  "
$ws.Range("D14").Value = "'25.89"
$ws.Range("E14").Value = "  -5.57%  "
$ws.Range("D15").Value = "'0.0000157"
$ws.Range("E15").Value = "  -10.18%  "
$ws.Range("D16").Value = "56.497.15"
$ws.Range("E16").Value = "  -2.90%  "
$ws.Range("E17").Value = "  -4.28%  "
$ws.Range("D18").Value = "2.944.71"
$ws.Range("E18").Value = "  -4.15%  "
$ws.Range("D19").Value = "'12.44"
$ws.Range("E19").Value = "  -5.74%  "
$ws.Range("E20").Value = "  -6.18%  "
$ws.Range("D21").Value = "'316.21"
$ws.Range("E21").Value = "  -7.81%  "
$ws.Range("E22").Value = "  -0.09%  "
$ws.Range("D23").Value = "'5.76"
$ws.Range("E23").Value = "  +0.16%  "
$ws.Range("E24").Value = "  -4.81%  "
$ws.Range("E25").Value = "  -4.83%  "
$ws.Range("E26").Value = "  +0.31%  "
$ws.Range("D27").Value = "'0.161"
$ws.Range("E27").Value = "  -5.26%  "
$ws.Range("D28").Value = "0.0₃0853"
$ws.Range("E28").Value = "  -12.88%  "
$ws.Range("E29").Value = "  -8.63%  "
$ws.Range("D30").Value = "'6.99"
$ws.Range("E30").Value = "  -6.82%  "
$ws.Range("D31").Value = "'1.74"
$ws.Range("E31").Value = "  -7.13%  "
$ws.Range("D32").Value = "'19.89"
$ws.Range("E32").Value = "  -6.34%  "
$ws.Range("D33").Value = "'1.13"
$ws.Range("E33").Value = "  -8.94%  "
$ws.Range("D34").Value = "'151.08"
$ws.Range("E34").Value = "  -4.46%  "
$ws.Range("D35").Value = "'4.43"
$ws.Range("E35").Value = "  -8.44%  "
$ws.Range("E36").Value = "  -5.97%  "
$ws.Range("E37").Value = "  -10.12%  "
$ws.Range("D38").Value = "'23.64"
$ws.Range("E38").Value = "  -10.31%  "
$ws.Range("E39").Value = "  -7.28%  "
$ws.Range("D40").Value = "'37.34"
$ws.Range("E40").Value = "  -1.40%  "
$ws.Range("D41").Value = "2.969.73"
$ws.Range("E41").Value = "  -4.44%  "
$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "  +0.11%  "
$ws.Range("D43").Value = "'3.67"
$ws.Range("E43").Value = "  -8.00%  "
$ws.Range("D44").Value = "'0.635"
$ws.Range("E44").Value = "  -5.39%  "
$ws.Range("D45").Value = "2.134.38"
$ws.Range("E45").Value = "  -8.77%  "
$ws.Range("D46").Value = "'1.34"
$ws.Range("E46").Value = "  -10.03%  "
$ws.Range("D47").Value = "'5.83"
$ws.Range("E47").Value = "  -4.35%  "
$ws.Range("D48").Value = "'0.906"
$ws.Range("E48").Value = "  -13.60%  "
$ws.Range("E49").Value = "  -6.10%  "
$ws.Range("D50").Value = "'18.82"
$ws.Range("E50").Value = "  -7.17%  "
$ws.Range("E51").Value = "  -6.64%  "
